$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "JSONMessageType"
$ws.Range("H2").Value = "DemoEvent"
$ws.Range("J6").Value = "JSONMessageType"
$ws.Range("H6").Value = "DemoEvent"
$ws.Range("J7").Value = "JSONMessageType"
$ws.Range("H7").Value = "DemoEvent"
